# Auto-generated script applying market price / profit updates across leve-profit sheets
# (mirrors a scheduled runner refreshing currentAveragePrice* / LeveProfit* columns)
$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "H10" = 900
    "I10" = 1000
    "J10" = 800
    "K10" = 1000
    "L10" = 800
    "M10" = -707
    "N10" = -1386
    "H19" = 1522
    "I19" = 1344.1428
    "J19" = 1699.8572
    "K19" = 1344.1428
    "L19" = 1699.8572
    "M19" = -1169.1428
    "N19" = -2049.8572
    "H33" = 213.25
    "I33" = 176
    "J33" = 399.5
    "K33" = 176
    "L33" = 399.5
    "M33" = 53
    "N33" = -857.5
    "H58" = 397.2
    "I58" = 121.75
    "J58" = 1499
    "K58" = 365.25
    "L58" = 4497
    "M58" = -215.25
    "N58" = -4797
    "I98" = 1115.3667
    "J98" = 0
    "K98" = 1115.3667
    "L98" = 0
    "M98" = 382.6333
    "H116" = 5001
    "I116" = 5001
    "J116" = 0
    "K116" = 5001
    "L116" = 0
    "M116" = -1559
    "I122" = 1115.3667
    "J122" = 0
    "K122" = 3346.1001
    "L122" = 0
    "M122" = -896.1001000000001
    "H132" = 3349.476
    "I132" = 3463.3333
    "J132" = 2666.3333
    "K132" = 10389.9999
    "L132" = 7998.999899999999
    "M132" = -7859.999899999999
    "N132" = -13058.9999
    "H138" = 4598.8
    "I138" = 1924.5652
    "J138" = 7394.591
    "K138" = 5773.6956
    "L138" = 22183.773
    "M138" = -633.6956
    "N138" = -32463.773
    "H141" = 2532.9285
    "I141" = 2532.9285
    "J141" = 0
    "K141" = 7598.7855
    "L141" = 0
    "M141" = -2418.7855
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
$deletedCells = @("N98", "N122")
foreach ($cellRef in $deletedCells) {
    $ws.Range($cellRef).ClearContents()
}

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "H24" = 38999.5
    "I24" = 0
    "J24" = 38999.5
    "K24" = 0
    "L24" = 38999.5
    "N24" = -39747.5
    "H32" = 1956.579
    "I32" = 1956.579
    "J32" = 0
    "K32" = 1956.579
    "L32" = 0
    "M32" = -1669.579
    "H97" = 149.375
    "I97" = 140.83333
    "J97" = 175
    "K97" = 140.83333
    "L97" = 175
    "M97" = 355.16667
    "N97" = -1167
    "H100" = 38999.5
    "I100" = 0
    "J100" = 38999.5
    "K100" = 0
    "L100" = 38999.5
    "N100" = -41163.5
    "H122" = 9512.125
    "I122" = 7299.5713
    "J122" = 25000
    "K122" = 21898.7139
    "L122" = 75000
    "M122" = -19448.7139
    "N122" = -79900
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "H22" = 4332067.5
    "I22" = 4940.4
    "J22" = 7938007
    "K22" = 4940.4
    "L22" = 7938007
    "M22" = -4767.4
    "N22" = -7938353
    "H94" = 6857.4443
    "I94" = 8459.846
    "J94" = 2691.2
    "K94" = 8459.846
    "L94" = 2691.2
    "M94" = -8008.846
    "N94" = -3593.2
    "H105" = 2889.5715
    "I105" = 2269.25
    "J105" = 3716.6667
    "K105" = 2269.25
    "L105" = 3716.6667
    "M105" = -522.25
    "N105" = -7210.6667
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "H31" = 4847.8335
    "I31" = 2869.5
    "J31" = 10782.833
    "K31" = 2869.5
    "L31" = 10782.833
    "M31" = -2574.5
    "N31" = -11372.833
    "H34" = 4847.8335
    "I34" = 2869.5
    "J34" = 10782.833
    "K34" = 2869.5
    "L34" = 10782.833
    "M34" = -2667.5
    "N34" = -11186.833
    "H86" = 13104.583
    "I86" = 9585.666999999999
    "J86" = 14277.556
    "K86" = 9585.666999999999
    "L86" = 14277.556
    "M86" = -8462.666999999999
    "N86" = -16523.556
    "H89" = 13104.583
    "I89" = 9585.666999999999
    "K89" = 47928.335
    "M89" = -42312.335
    "H94" = 2081.4167
    "I94" = 2279.4
    "J94" = 1940
    "K94" = 2279.4
    "L94" = 1940
    "M94" = -1828.4
    "N94" = -2842
    "H134" = 14708587
    "I134" = 17859716
    "J134" = 3314.3333
    "K134" = 53579148
    "L134" = 9942.999899999999
    "M134" = -53576613
    "N134" = -15012.9999
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$updates = @{
    "H3" = 10116.429
    "I3" = 9302.5
    "J3" = 15000
    "K3" = 27907.5
    "L3" = 45000
    "M3" = -27795.5
    "N3" = -45224
    "H5" = 101192.5
    "I5" = 200429.4
    "J5" = 1955.6
    "K5" = 601288.2
    "L5" = 5866.799999999999
    "M5" = -601176.2
    "N5" = -6090.799999999999
    "H106" = 23999.445
    "I106" = 0
    "J106" = 23999.445
    "K106" = 0
    "L106" = 71998.33499999999
    "N106" = -73890.33499999999
    "H135" = 101192.5
    "I135" = 200429.4
    "J135" = 1955.6
    "K135" = 1803864.6
    "L135" = 17600.4
    "M135" = -1801329.6
    "N135" = -22670.4
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "H63" = 0
    "I63" = 0
    "J63" = 0
    "K63" = 0
    "L63" = 0
    "H66" = 0
    "I66" = 0
    "J66" = 0
    "K66" = 0
    "L66" = 0
    "H80" = 3204.182
    "I80" = 3412.125
    "J80" = 2649.6667
    "K80" = 3412.125
    "L80" = 2649.6667
    "M80" = -2414.125
    "N80" = -4645.6667
    "H83" = 3204.182
    "I83" = 3412.125
    "J83" = 2649.6667
    "K83" = 17060.625
    "L83" = 13248.3335
    "M83" = -12068.625
    "N83" = -23232.3335
    "H97" = 1894
    "I97" = 1688.6
    "J97" = 4975
    "K97" = 1688.6
    "L97" = 4975
    "M97" = -1192.6
    "N97" = -5967
    "H132" = 2981444.2
    "I132" = 3294899.2
    "J132" = 3621.5
    "K132" = 9884697.600000001
    "L132" = 10864.5
    "M132" = -9882167.600000001
    "N132" = -15924.5
    "H139" = 130337.4
    "I139" = 0
    "J139" = 130337.4
    "K139" = 0
    "L139" = 130337.4
    "N139" = -140617.4
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
$deletedCells = @("N63", "N66")
foreach ($cellRef in $deletedCells) {
    $ws.Range($cellRef).ClearContents()
}

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "H16" = 2727.1667
    "I16" = 1590.75
    "J16" = 5000
    "K16" = 1590.75
    "L16" = 5000
    "M16" = -1420.75
    "N16" = -5340
    "H46" = 3731.25
    "I46" = 2721.875
    "J46" = 5750
    "K46" = 2721.875
    "L46" = 5750
    "M46" = -2533.875
    "N46" = -6126
    "H76" = 19064.75
    "I76" = 7760
    "J76" = 22833
    "K76" = 7760
    "L76" = 22833
    "M76" = -7422
    "N76" = -23509
    "H79" = 19064.75
    "I79" = 7760
    "J79" = 22833
    "K79" = 7760
    "L79" = 22833
    "M79" = -6590
    "N79" = -25173
    "H132" = 41687204
    "I132" = 41687204
    "J132" = 0
    "K132" = 125061612
    "L132" = 0
    "M132" = -125059082
    "H133" = 89999
    "I133" = 0
    "J133" = 89999
    "K133" = 0
    "L133" = 89999
    "N133" = -95059
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "H126" = 4091.077
    "I126" = 4048.6667
    "J126" = 4600
    "K126" = 12146.0001
    "L126" = 13800
    "M126" = -9676.000100000001
    "N126" = -18740
    "H136" = 10418342
    "I136" = 10871283
    "J136" = 694.5
    "K136" = 32613849
    "L136" = 2083.5
    "M136" = -32611299
    "N136" = -7183.5
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
